# Automatic update of files.
#
# Upstream, the underlying export re-paired several observation rows in
# the "Artfynd" sheet: the row *positions* (row numbers) are unchanged,
# but the full data of certain rows now belongs to what used to be a
# different row. Reproduce that by swapping / rotating the complete row
# contents (columns A:AY) between the affected rows.
#
# Columns Y ("Startdatum") and AA ("Slutdatum") hold a plain-text date
# string (e.g. "2026-01-24"), not a real Excel date. A straight
# Range.Value2 assignment lets Excel "smart type" that string into a
# real date serial, which would corrupt the cell type, so those two
# columns are re-applied as text explicitly after each row copy.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Get-RowValues($rowNum) {
    return $ws.Range("A" + $rowNum + ":AY" + $rowNum).Value2
}

function Set-RowValues($rowNum, $values, $yText, $aaText) {
    $ws.Range("A" + $rowNum + ":AY" + $rowNum).Value2 = $values

    # Re-assert Y/AA as literal text so Excel doesn't reinterpret the
    # "yyyy-mm-dd" looking string as a real date value.
    $yCell = $ws.Range("Y" + $rowNum)
    $yCell.NumberFormat = "@"
    $yCell.Value2 = $yText
    $yCell.NumberFormat = "General"

    $aaCell = $ws.Range("AA" + $rowNum)
    $aaCell.NumberFormat = "@"
    $aaCell.Value2 = $aaText
    $aaCell.NumberFormat = "General"
}

# Capture full row snapshots (values + the Y/AA text) before any writes,
# since rows are paired/rotated (not every row maps 1:1 to its immediate
# neighbor).
$snapshot = @{}
foreach ($r in 13, 14, 20, 21, 22, 23, 24, 27, 28) {
    $snapshot[$r] = @{
        Values = Get-RowValues $r
        Y      = $ws.Range("Y" + $r).Value2
        AA     = $ws.Range("AA" + $r).Value2
    }
}

# --- Rows 13 <-> 14 : simple swap ---
Set-RowValues 13 $snapshot[14].Values $snapshot[14].Y $snapshot[14].AA
Set-RowValues 14 $snapshot[13].Values $snapshot[13].Y $snapshot[13].AA

# --- Rows 20 <-> 21 : simple swap ---
Set-RowValues 20 $snapshot[21].Values $snapshot[21].Y $snapshot[21].AA
Set-RowValues 21 $snapshot[20].Values $snapshot[20].Y $snapshot[20].AA

# --- Rows 22, 23, 24 : 3-way rotation ---
# new row22 = old row24, new row23 = old row22, new row24 = old row23
Set-RowValues 22 $snapshot[24].Values $snapshot[24].Y $snapshot[24].AA
Set-RowValues 23 $snapshot[22].Values $snapshot[22].Y $snapshot[22].AA
Set-RowValues 24 $snapshot[23].Values $snapshot[23].Y $snapshot[23].AA

# --- Rows 27 <-> 28 : simple swap ---
Set-RowValues 27 $snapshot[28].Values $snapshot[28].Y $snapshot[28].AA
Set-RowValues 28 $snapshot[27].Values $snapshot[27].Y $snapshot[27].AA
